# A new price-record row is inserted at row 153 of the sheet (pushing the
# existing rows 153:227 down to 154:228, and extending the used range from
# A1:R227 to A1:R228). The new row carries a fresh "Cultivar IV Región /
# Primera" quote dated 2022-05-23 (serial 44704).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 153; everything below shifts down one row.
$ws.Rows("153:153").Insert()

# Populate the newly inserted row 153 with the new record.
$ws.Cells.Item(153, 1).Value2  = 10
$ws.Cells.Item(153, 2).Value2  = "Vega Modelo de Temuco"
$ws.Cells.Item(153, 3).Value2  = "La Araucanía"
$ws.Cells.Item(153, 4).Value2  = 44704
$ws.Cells.Item(153, 5).Value2  = 9
$ws.Cells.Item(153, 6).Value2  = 100112043
$ws.Cells.Item(153, 7).Value2  = "Pepino dulce"
$ws.Cells.Item(153, 8).Value2  = "Cultivar IV Región"
$ws.Cells.Item(153, 9).Value2  = "Primera"
$ws.Cells.Item(153, 10).Value2 = 200
$ws.Cells.Item(153, 11).Value2 = 19000
$ws.Cells.Item(153, 12).Value2 = 19000
$ws.Cells.Item(153, 13).Value2 = 19000
$ws.Cells.Item(153, 14).Value2 = "$/bandeja 18 kilos"
$ws.Cells.Item(153, 15).Value2 = "Provincia de Limarí"
$ws.Cells.Item(153, 16).Value2 = 1056
$ws.Cells.Item(153, 17).Value2 = 18
$ws.Cells.Item(153, 18).Value2 = "Hortaliza"
